$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44314
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806

# Row 4
$ws.Range("D4").Value = 44229
$ws.Range("K4").Value = "Fortuna"

# Row 5
$ws.Range("D5").Value = 44245
$ws.Range("K5").Value = "Black Amber"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 250

# Row 7
$ws.Range("D7").Value = 44243
$ws.Range("K7").Value = "Black Amber"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 14500
$ws.Range("S7").Value = 806

# Row 8
$ws.Range("D8").Value = 44174
$ws.Range("K8").Value = "Angeleno"
$ws.Range("M8").Value = 270
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("S8").Value = 1139

# Row 9
$ws.Range("D9").Value = 44285
$ws.Range("K9").Value = "Angeleno"
$ws.Range("Q9").Value = "$/bandeja 18 kilos granel"

# Row 10
$ws.Range("D10").Value = 44278
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = "$/caja 18 kilos granel"
$ws.Range("S10").Value = 861

# Row 11
$ws.Range("D11").Value = 44238
$ws.Range("K11").Value = "Black Amber"
$ws.Range("L11").Value = "Segunda"
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14500
$ws.Range("Q11").Value = "$/bandeja 18 kilos granel"
$ws.Range("S11").Value = 806

# Row 12
$ws.Range("D12").Value = 44238

# Row 13
$ws.Range("D13").Value = 44175
$ws.Range("K13").Value = "Angeleno"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("S13").Value = 1194

# Row 14
$ws.Range("D14").Value = 44239
$ws.Range("K14").Value = "Fortuna"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 15500
$ws.Range("S14").Value = 861
